$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3.08
    3  = 3.21
    4  = 3.14
    5  = 3.17
    6  = 3.09
    7  = 3.26
    8  = 3.17
    9  = 2.96
    10 = 3.14
    11 = 3.11
    12 = 3.15
    13 = 3.05
    14 = 3
    15 = 3.09
    16 = 3.13
    17 = 3.06
    18 = 2.97
    19 = 3.1
    20 = 1.84
    21 = 2.3
    22 = 2.22
    23 = 2.27
    24 = 2.25
    25 = 2.29
    26 = 6.6
    27 = 3.49
    28 = 3.79
    29 = 2.79
    30 = 3.66
    31 = 8.5
    32 = 2.51
    33 = 2.56
    34 = 2.98
    35 = 3.75
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
